$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update a few subcategory (column H) labels
$ws.Range("H15").Value = "mixed statistical plot (more than 1 statistical plot and type)"
$ws.Range("H16").Value = "data display"
$ws.Range("H21").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H25").Value = "photo(s)"
$ws.Range("H29").Value = "line graph(s)"
$ws.Range("H30").Value = "data display"
$ws.Range("H34").Value = "photo(s)"

# Remove the "is_viewed" column (column I) entirely
$ws.Range("I1:I34").EntireColumn.Delete()
